$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D" = 0.105
    "G" = 0.01374045801526718
    "H" = 0.01374045801526718
    "I" = -0.001475826972010178
    "J" = -0.001475826972010178
    "K" = -2.77
    "L" = -0.03524173027989822
    "M" = 0.024
    "N" = 0.002142857142857143
    "O" = -0.008664259927797834
    "S" = 0.024
    "T" = 1
    "U" = 7.14
    "V" = 0.6375000000000001
    "W" = -0.4973070017953321
    "X" = 0.1850293811415115
    "Y" = -0.6823363829368436
    "Z" = 2.815186246418338
    "AA" = -0.004154727793696275
    "AB" = 0.08233514316679655
    "AC" = -0.08648987096049282
    "AD" = 26.4
    "AF" = 26.4
    "AG" = 19.26
    "AH" = 0.7021276595744681
    "AI" = 0.9056603773584906
    "AJ" = 0.6323046618516086
    "AK" = 0.8750567923671059
    "AL" = 1.05
    "AM" = 0.9730000000000001
    "AN" = 28.32618025751073
    "AO" = -0.1104761904761905
    "AP" = 20.66523605150214
    "AQ" = -0.1192189105858171
}

foreach ($row in 2, 3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
